$d = $word.ActiveDocument

# The document currently starts with:
#   Paragraph 1 (Heading1): "On Pilgrimage - March/April 1978"
#   Paragraph 2 (bold run):  "By Dorothy Day"
#
# We replace both paragraphs in one shot with the pandoc-style title block:
#   Paragraph 1 (Title style):   "On" " " "Pilgrimage" " " "-" " " "March" "/" "April" " " "1978"
#     -- each token its own run, mirroring the target markup produced by the
#        docx writer that generates one run per token.
#   Paragraph 2 (Authors style): "Dorothy" " " "Day" (no "By " prefix, no bold)

$wNs = "xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'"

function Run([string]$text) {
    return "<w:r><w:t xml:space='preserve'>$text</w:t></w:r>"
}

$titleRuns = (Run "On") + (Run " ") + (Run "Pilgrimage") + (Run " ") + (Run "-") + (Run " ") + `
             (Run "March") + (Run "/") + (Run "April") + (Run " ") + (Run "1978")

$authorRuns = (Run "Dorothy") + (Run " ") + (Run "Day")

$titlePara  = "<w:p $wNs><w:pPr><w:pStyle w:val='Title'/></w:pPr>$titleRuns</w:p>"
$authorPara = "<w:p $wNs><w:pPr><w:pStyle w:val='Authors'/></w:pPr>$authorRuns</w:p>"

$titleP = $d.Paragraphs(1)
$authorP = $d.Paragraphs(2)

$full = $d.Range($titleP.Range.Start, $authorP.Range.End)
$full.InsertXML($titlePara + $authorPara)
